$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.346.17"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "2.493.29"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.94%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").Value = "2.493.02"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.70%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("E13").Value = "  -3.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "2.947.11"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").Value = "67.122.27"
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "2.492.14"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("E23").Value = "  -6.60%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("E30").Value = "  -6.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "505.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.04%  "
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.58"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.337"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.81%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.36%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.540"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.54%  "
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("E51").Value = "  -4.09%  "
